# Update the "想去人数" (interested-count) values in column F for the
# events that changed between scrapes, on both the "展览" (sheet1) and
# "全部类型" (sheet4) sheets — they carry duplicate rows for the same
# events, so both need the same bump.

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览"   = @{ "F3" = 210; "F4" = 694; "F7" = 2840; "F9" = 8100; "F10" = 209; "F11" = 479; "F13" = 421 }
    "全部类型" = @{ "F3" = 210; "F4" = 694; "F9" = 2840; "F11" = 8100; "F12" = 209; "F13" = 479; "F17" = 421 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellUpdates = $updates[$sheetName]
    foreach ($cellRef in $cellUpdates.Keys) {
        $ws.Range($cellRef).Value = $cellUpdates[$cellRef]
    }
}
